$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Key Metrics": add a new "Q4FY22" column (D) that duplicates the
# existing "Q4FY23" column (C), mirroring the new KPI quarter being tested.
# ---------------------------------------------------------------------------
$wsKM = $wb.Worksheets.Item("Key Metrics")

$wsKM.Range("D1").Value = "Q4FY22"

$kmRows = 2,3,4,5,6,7,8,9,10,11,12
foreach ($r in $kmRows) {
    $src = $wsKM.Cells.Item($r, 3)
    $dst = $wsKM.Cells.Item($r, 4)
    $dst.Value = $src.Value2
    if ($src.NumberFormat -ne "General") {
        $dst.NumberFormat = $src.NumberFormat
    }
}

[void]$wsKM.Range("B1:D1").Select()

# ---------------------------------------------------------------------------
# Sheet "Balance Sheet": add new quarterly columns E (Q4FY24), F (Q4FY23),
# G (Q4FY22) to both tables (rows 1-8 and 11-18), copying from the existing
# annual columns B (FY24) and C (FY23) as placeholder KPI data.
# ---------------------------------------------------------------------------
$wsBS = $wb.Worksheets.Item("Balance Sheet")

foreach ($r in 1, 11) {
    $h1 = $wsBS.Cells.Item($r, 5)
    $h1.Value = "Q4FY24"
    $h1.Font.Bold = $true
    $h1.HorizontalAlignment = -4108

    $h2 = $wsBS.Cells.Item($r, 6)
    $h2.Value = "Q4FY23"
    $h2.Font.Bold = $true
    $h2.HorizontalAlignment = -4108

    $h3 = $wsBS.Cells.Item($r, 7)
    $h3.Value = "Q4FY22"
    $h3.Font.Bold = $true
    $h3.HorizontalAlignment = -4108
}

$bsDataRows = 2,3,4,5,6,7,8,12,13,14,15,16,17,18
foreach ($r in $bsDataRows) {
    $srcB = $wsBS.Cells.Item($r, 2)
    $srcC = $wsBS.Cells.Item($r, 3)
    $dstE = $wsBS.Cells.Item($r, 5)
    $dstF = $wsBS.Cells.Item($r, 6)
    $dstG = $wsBS.Cells.Item($r, 7)

    $dstE.Value = $srcB.Value2
    if ($srcB.NumberFormat -ne "General") { $dstE.NumberFormat = $srcB.NumberFormat }

    $dstF.Value = $srcC.Value2
    if ($srcC.NumberFormat -ne "General") { $dstF.NumberFormat = $srcC.NumberFormat }

    $dstG.Value = $srcC.Value2
    if ($srcC.NumberFormat -ne "General") { $dstG.NumberFormat = $srcC.NumberFormat }
}

[void]$wsBS.Range("E11:G11").Select()
